# Add a new row (33) of master data to the single worksheet, mirroring the
# pattern of the existing rows (regcntr_id, machine_id, lang_code, is_active,
# cr_by, cr_dtimes, eff_dtimes).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Matches the saved selection recorded in the workbook after the edit.
$ws.Range("B30").Select()
